# Purchase-request / vendor product template: add a hidden "helper" list
# (product_category options) in column CO and wire a data-validation
# dropdown on C2:C1500 to it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper data in column CO (93) ---------------------------------------
# F2 gets a default numeric value (total_quantity_purchased starter row)
$ws.Range("F2").Value = 0

# CO2 mirrors the "product_category" header (bold, like row 1 header cells)
$ws.Range("CO2").Value = "product_category"
$ws.Range("CO2").Font.Bold = $true

# The three valid category options living in CO3:CO5 feed the validation
# list below. Write CO3 then CO5 then CO4 so new shared-string entries land
# in the same order as the source workbook (consumable, service Product,
# stockable).
$ws.Range("CO3").Value = "consumable"
$ws.Range("CO5").Value = "service Product"
$ws.Range("CO4").Value = "stockable"

# Column CO should be sized like the other text columns.
$ws.Columns.Item(93).ColumnWidth = 15.7

# --- dropdown validation on the product_category column ------------------
$validation = $ws.Range("C2:C1500").Validation
$validation.Add(3, 1, 1, "=`$CO`$3:`$CO`$5")
$validation.ErrorTitle = "You have to select from list"
$validation.InputTitle = "select from list"

# --- leave the selection where the author left it -------------------------
$null = $ws.Range("C4").Select()
